$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the "Normal" style to the new row range first so the written cells
# pick up a dedicated cellXf (mirrors the style index used by the target
# workbook's new row).
$ws.Range("A4:I4").Style = "Normal"

# New row of import data (numeric-cell regression test fixture).
$ws.Range("A4").Value = "designer"
$ws.Range("B4").Value = 12345
$ws.Range("C4").Value = 12345
$ws.Range("D4").Value = "Magasin Ikea"
$ws.Range("E4").Value = "Test with numeric values"
$ws.Range("G4").Value = "blue"
$ws.Range("H4").Value = 12345
$ws.Range("I4").Value = "blue,green"

# Move the active selection like the author's workbook ended up with.
$ws.Range("E5").Select()
